$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 2 is the "https://pypi.org/project/webdriver-manager/" bullet.
# Insert a new bullet paragraph right after it with the new URL, inheriting
# the same paragraph/run formatting (bullet, font, size) as paragraph 2.
$para2 = $tr.Paragraphs(2, 1)
$newRange = $para2.InsertAfter("`rhttps://www.selenium.dev/selenium/docs/api/rb/Selenium/WebDriver/Chrome/Options.html")
